{"js": "// Remove the trailing \"The results ... because of it.\" text from the\n// Stage 2 paragraph, leaving \"...available separately as well. \" intact\n// (with its trailing space) and leaving the `_GoBack` bookmark in place.\n\nconst body = context.document.body;\n\n// The original text spans two runs and crosses the `_GoBack` bookmark:\n//   \"...available separately as well.\"\n//   + \" The results were not very satisfying, and if it worked for some\n//      pictures, in did not correct others, I included an example of the\n//      latt\" [_GoBack bookmark] \"er, to work on it and correct possible\n//      errors of the plugin, if the behavior is because of it.\"\n//\n// We delete the text in two passes so the bookmark (which sits between the\n// two halves of the word \"latter\") is left untouched: first remove the\n// portion that follows the bookmark, then the portion that precedes it.\n\nconst afterBookmark =\n  \"er, to work on it and correct possible errors of the plugin, \" +\n  \"if the behavior is because of it.\";\nconst afterResults = body.search(afterBookmark, { matchCase: true });\nafterResults.load(\"text\");\nawait context.sync();\nif (afterResults.items.length > 0) {\n  afterResults.items[0].delete();\n}\nawait context.sync();\n\nconst beforeBookmark =\n  \"The results were not very satisfying, and if it worked for some \" +\n  \"pictures, in did not correct others, I included an example of the latt\";\nconst beforeResults = body.search(beforeBookmark, { matchCase: true });\nbeforeResults.load(\"text\");\nawait context.sync();\nif (beforeResults.items.length > 0) {\n  beforeResults.items[0].delete();\n}\nawait context.sync();\n", "ps1": "# Remove the trailing \"The results ... because of it.\" text from the\n# Stage 2 paragraph, leaving \"...available separately as well. \" intact\n# (with its trailing space) and leaving the `_GoBack` bookmark in place.\n#\n# The original text spans two runs and crosses the `_GoBack` bookmark:\n#   \"...available separately as well.\"\n#   + \" The results were not very satisfying, and if it worked for some\n#      pictures, in did not correct others, I included an example of the\n#      latt\" [_GoBack bookmark] \"er, to work on it and correct possible\n#      errors of the plugin, if the behavior is because of it.\"\n#\n# We remove the text in two Find/Replace passes so the bookmark (which sits\n# between the two halves of the word \"latter\") is left untouched: first the\n# portion that follows the bookmark, then the portion that precedes it.\n\n$d = $word.ActiveDocument\n\n$rngAfter = $d.Content\n$rngAfter.Find.ClearFormatting()\n$rngAfter.Find.Replacement.ClearFormatting()\n$rngAfter.Find.Text = \"er, to work on it and correct possible errors of the plugin, if the behavior is because of it.\"\n$rngAfter.Find.Replacement.Text = \"\"\n$rngAfter.Find.Execute([ref]$rngAfter.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$rngAfter.Find.Replacement.Text, 2)\n\n$rngBefore = $d.Content\n$rngBefore.Find.ClearFormatting()\n$rngBefore.Find.Replacement.ClearFormatting()\n$rngBefore.Find.Text = \"The results were not very satisfying, and if it worked for some pictures, in did not correct others, I included an example of the latt\"\n$rngBefore.Find.Replacement.Text = \"\"\n$rngBefore.Find.Execute([ref]$rngBefore.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$rngBefore.Find.Replacement.Text, 2)\n"}
